$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Refresh the cryptos table: update prices / 1h volume %, and fix two
# rows (Toncoin/InjectiveProtocol and TheGraph/WEMIXToken) that had
# swapped data in the source feed.
# A leading apostrophe forces Excel to store the value as literal text
# (preserving formats like '1.80' or '51.080.17'); Style is reset to
# 'Normal' afterwards so no stray number-format/quote-prefix style sticks.
$ws.Range('D2').Value = "'51.126.23"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.77%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.946.27"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.21%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.04%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'376.98"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.39%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'101.55"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.60%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.539"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -1.22%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.06%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.592"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.62%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'36.46"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.32%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.54%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.0854"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.92%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.407.42"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.22%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'18.18"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -1.24%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'7.63"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.81%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.937.24"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.56%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.997"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.32%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'10.95"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +47.02%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'51.080.17"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.72%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -6.61%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'12.51"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -3.09%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.0₃0956"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.92%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'266.09"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.23%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'68.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.19%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'3.14"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +7.81%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -2.04%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'7.47"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -3.59%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.01%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.57%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'25.68"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -1.25%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -4.99%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.25%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'50.79"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.41%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'Toncoin"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'2.05"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.78%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'InjectiveProtocol"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'33.48"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -3.89%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  -1.98%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.20%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +4.33%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.15%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'16.45"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -3.70%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'1.80"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -2.50%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.49"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -3.83%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'120.26"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -1.54%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'21.35"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.86%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'3.39"
$ws.Range('D45').Style = 'Normal'
$ws.Range('B46').Value = "'WEMIXToken"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'2.03"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.24%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'TheGraph"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'0.273"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.66%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'2.31"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.71%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.991.90"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -2.02%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.46%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.91%  "
$ws.Range('E51').Style = 'Normal'
